$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2751.5
$ws.Range("I62").Value = 2402.8
$ws.Range("J62").Value = 3332.6667
$ws.Range("K62").Value = 2402.8
$ws.Range("L62").Value = 3332.6667
$ws.Range("M62").Value = -1778.8
$ws.Range("N62").Value = -4580.6667
$ws.Range("H65").Value = 2751.5
$ws.Range("I65").Value = 2402.8
$ws.Range("J65").Value = 3332.6667
$ws.Range("K65").Value = 12014
$ws.Range("L65").Value = 16663.3335
$ws.Range("M65").Value = -8894
$ws.Range("N65").Value = -22903.3335
$ws.Range("H92").Value = 494.57144
$ws.Range("I92").Value = 449
$ws.Range("J92").Value = 640.4
$ws.Range("K92").Value = 449
$ws.Range("L92").Value = 640.4
$ws.Range("M92").Value = 799
$ws.Range("N92").Value = -3136.4
$ws.Range("H112").Value = 3998.5715
$ws.Range("J112").Value = 4098.3335
$ws.Range("L112").Value = 12295.0005
$ws.Range("N112").Value = -14511.0005
$ws.Range("H116").Value = 54832132
$ws.Range("I116").Value = 60195184
$ws.Range("K116").Value = 60195184
$ws.Range("M116").Value = -60191742
$ws.Range("H135").Value = 3412.4
$ws.Range("I135").Value = 928.36365
$ws.Range("J135").Value = 10243.5
$ws.Range("K135").Value = 8355.272849999999
$ws.Range("L135").Value = 92191.5
$ws.Range("M135").Value = -5820.272849999999
$ws.Range("N135").Value = -97261.5
$ws.Range("H136").Value = 114352
$ws.Range("J136").Value = 114352
$ws.Range("L136").Value = 114352
$ws.Range("N136").Value = -124552
$ws.Range("H139").Value = 123799.8
$ws.Range("J139").Value = 123799.8
$ws.Range("L139").Value = 123799.8
$ws.Range("N139").Value = -134079.8
$ws.Range("H140").Value = 56852.125
$ws.Range("J140").Value = 54872.57
$ws.Range("L140").Value = 54872.57
$ws.Range("N140").Value = -65232.57

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1802334.8
$ws.Range("I2").Value = 2499906.8
$ws.Range("K2").Value = 2499906.8
$ws.Range("M2").Value = -2499793.8
$ws.Range("H74").Value = 1200.5416
$ws.Range("I74").Value = 962.06665
$ws.Range("J74").Value = 1598
$ws.Range("K74").Value = 962.06665
$ws.Range("L74").Value = 1598
$ws.Range("M74").Value = -88.06664999999998
$ws.Range("N74").Value = -3346
$ws.Range("H77").Value = 1200.5416
$ws.Range("I77").Value = 962.06665
$ws.Range("J77").Value = 1598
$ws.Range("K77").Value = 4810.33325
$ws.Range("L77").Value = 7990
$ws.Range("M77").Value = -442.3332499999997
$ws.Range("N77").Value = -16726
$ws.Range("H109").Value = 68333
$ws.Range("J109").Value = 68333
$ws.Range("L109").Value = 68333
$ws.Range("N109").Value = -71107
$ws.Range("H116").Value = 1802334.8
$ws.Range("I116").Value = 2499906.8
$ws.Range("K116").Value = 2499906.8
$ws.Range("M116").Value = -2497612.8
$ws.Range("H122").Value = 4063.0625
$ws.Range("I122").Value = 2111.9546
$ws.Range("J122").Value = 8355.5
$ws.Range("K122").Value = 6335.8638
$ws.Range("L122").Value = 25066.5
$ws.Range("M122").Value = -3885.8638
$ws.Range("N122").Value = -29966.5
$ws.Range("H140").Value = 114250
$ws.Range("J140").Value = 114250
$ws.Range("L140").Value = 114250
$ws.Range("N140").Value = -124610

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1802334.8
$ws.Range("I3").Value = 2499906.8
$ws.Range("K3").Value = 2499906.8
$ws.Range("M3").Value = -2499792.8
$ws.Range("H86").Value = 1270
$ws.Range("I86").Value = 1450
$ws.Range("J86").Value = 1250
$ws.Range("K86").Value = 1450
$ws.Range("L86").Value = 1250
$ws.Range("M86").Value = -327
$ws.Range("N86").Value = -3496
$ws.Range("H89").Value = 1270
$ws.Range("I89").Value = 1450
$ws.Range("J89").Value = 1250
$ws.Range("K89").Value = 7250
$ws.Range("L89").Value = 6250
$ws.Range("M89").Value = -1634
$ws.Range("N89").Value = -17482
$ws.Range("H94").Value = 392102.7
$ws.Range("I94").Value = 913734.4
$ws.Range("K94").Value = 913734.4
$ws.Range("M94").Value = -913283.4
$ws.Range("H107").Value = 1284.1892
$ws.Range("J107").Value = 1246.5
$ws.Range("L107").Value = 1246.5
$ws.Range("N107").Value = -5086.5
$ws.Range("H140").Value = 98147.836
$ws.Range("J140").Value = 98147.836
$ws.Range("L140").Value = 98147.836
$ws.Range("N140").Value = -108507.836

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 91044
$ws.Range("I141").Value = 39498
$ws.Range("J141").Value = 94725.86
$ws.Range("K141").Value = 39498
$ws.Range("L141").Value = 94725.86
$ws.Range("M141").Value = -34318
$ws.Range("N141").Value = -105085.86

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 5000
$ws.Range("J17").Value = 5000
$ws.Range("L17").Value = 15000
$ws.Range("N17").Value = -15338
$ws.Range("H34").Value = 1964.3334
$ws.Range("J34").Value = 5000
$ws.Range("L34").Value = 15000
$ws.Range("N34").Value = -15168
$ws.Range("H55").Value = 2000
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").Value = $null
$ws.Range("H95").Value = 16666
$ws.Range("J95").Value = 17999.2
$ws.Range("L95").Value = 53997.60000000001
$ws.Range("N95").Value = -58115.60000000001
$ws.Range("H126").Value = 24999.666
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 24999.666
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 74998.99800000001
$ws.Range("M126").Value = $null
$ws.Range("N126").Value = -84878.99800000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 2511729
$ws.Range("I70").Value = 3180023.5
$ws.Range("J70").Value = 5623.75
$ws.Range("K70").Value = 3180023.5
$ws.Range("L70").Value = 5623.75
$ws.Range("M70").Value = -3179753.5
$ws.Range("N70").Value = -6163.75
$ws.Range("H73").Value = 2511729
$ws.Range("I73").Value = 3180023.5
$ws.Range("J73").Value = 5623.75
$ws.Range("K73").Value = 3180023.5
$ws.Range("L73").Value = 5623.75
$ws.Range("M73").Value = -3179087.5
$ws.Range("N73").Value = -7495.75
$ws.Range("H92").Value = 19999.5
$ws.Range("J92").Value = 19999.5
$ws.Range("L92").Value = 19999.5
$ws.Range("N92").Value = -23743.5
$ws.Range("H102").Value = 4698.3335
$ws.Range("I102").Value = 4754.3667
$ws.Range("K102").Value = 4754.3667
$ws.Range("M102").Value = -3132.3667
$ws.Range("H113").Value = 12374.75
$ws.Range("I113").Value = 9750
$ws.Range("J113").Value = 14999.5
$ws.Range("K113").Value = 9750
$ws.Range("L113").Value = 14999.5
$ws.Range("M113").Value = -7580
$ws.Range("N113").Value = -19339.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5724.2188
$ws.Range("I46").Value = 3200.2856
$ws.Range("J46").Value = 6430.92
$ws.Range("K46").Value = 3200.2856
$ws.Range("L46").Value = 6430.92
$ws.Range("M46").Value = -3012.2856
$ws.Range("N46").Value = -6806.92
$ws.Range("H122").Value = 95242024
$ws.Range("I122").Value = 125003150
$ws.Range("K122").Value = 375009450
$ws.Range("M122").Value = -375007000
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2055.3125
$ws.Range("I107").Value = 2331.52
$ws.Range("K107").Value = 6994.559999999999
$ws.Range("M107").Value = -5074.559999999999
$ws.Range("H122").Value = 3726.2222
$ws.Range("I122").Value = 3075.55
$ws.Range("K122").Value = 9226.650000000001
$ws.Range("M122").Value = -6776.650000000001
$ws.Range("H126").Value = 1457.4286
$ws.Range("I126").Value = 1200.3334
$ws.Range("K126").Value = 3601.0002
$ws.Range("M126").Value = -1131.0002
$ws.Range("H132").Value = 15155312
$ws.Range("I132").Value = 1643.2667
$ws.Range("K132").Value = 4929.800099999999
$ws.Range("M132").Value = -2399.800099999999
